$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config_BC")
$ws.Activate()
$excel.Goto($ws.Range("C1:C80"))
